$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The existing data (rows 2-7) is being extended with a fresh batch of
# Samity rows (8-13) that reuse the same Branch/Samity codes, picking the
# "System Generated Samity Information" values back up where the last
# batch (222-007..222-010) left off: 222-011..222-014 for the four
# migrated rows, and the last two rows duplicate the "not migrated"
# rows (6-7) exactly.

# Copy the whole A2:D7 block down to A8:D13 first so columns A, B and D
# (which are unchanged) keep their original shared-string text values.
$ws.Range("A2:D7").Copy()
$ws.Range("A8").PasteSpecial()

# Now patch up column C ("System Generated Samity Information") for the
# four newly migrated rows.
$ws.Cells.Item(8, 3).Value = "222-011 - LAMP 1"
$ws.Cells.Item(9, 3).Value = "222-012 - Dollan Chapa"
$ws.Cells.Item(10, 3).Value = "222-013 - Golap"
$ws.Cells.Item(11, 3).Value = "222-014 - Modhomoti"

# Rows 12 and 13 mirror rows 6 and 7 (not migrated yet): column C there
# stays blank text, same as the source rows. A bare "'" forces an empty
# text value instead of clearing the cell outright; ClearFormats drops
# the transient quote-prefix formatting so the cell keeps the default
# style.
$ws.Cells.Item(12, 3).Value = "'"
$ws.Cells.Item(12, 3).ClearFormats()
$ws.Cells.Item(13, 3).Value = "'"
$ws.Cells.Item(13, 3).ClearFormats()
